# Apply updated "Tp2" averages across the three summary sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Promedio_edad ---
$ws1 = $wb.Worksheets.Item("Promedio_edad")

$ws1.Range("B2").Value = 116
$ws1.Range("C2").Value = 4.224137931034483
$ws1.Range("D2").Value = 4.22

$ws1.Range("C3").Value = 4.90940170940171
$ws1.Range("D3").Value = 4.91

$ws1.Range("C4").Value = 5.014553014553014
$ws1.Range("D4").Value = 5.01

$ws1.Range("C5").Value = 4.717241379310344
$ws1.Range("D5").Value = 4.72

$ws1.Range("C6").Value = 4.434782608695652
$ws1.Range("D6").Value = 4.43

# --- Sheet 2: Promedio_genero ---
$ws2 = $wb.Worksheets.Item("Promedio_genero")

$ws2.Range("B2").Value = 969
$ws2.Range("C2").Value = 4.863777089783282
$ws2.Range("D2").Value = 4.86

$ws2.Range("B3").Value = 740
$ws2.Range("C3").Value = 4.758108108108108
$ws2.Range("D3").Value = 4.76

# --- Sheet 3: Promedio_ocup ---
$ws3 = $wb.Worksheets.Item("Promedio_ocup")

$ws3.Range("C2").Value = 5.014218009478673
$ws3.Range("D2").Value = 5.01

$ws3.Range("C3").Value = 4.893817204301075
$ws3.Range("D3").Value = 4.89

$ws3.Range("C4").Value = 4.647540983606557
$ws3.Range("D4").Value = 4.65

# Row 5 and Row 6 swap: D.Estudiante now sits on row 5, E.Jubilado moves to row 6
$ws3.Range("A5").Value = "D.Estudiante"
$ws3.Range("B5").Value = 272
$ws3.Range("C5").Value = 4.544117647058823
$ws3.Range("D5").Value = 4.54

$ws3.Range("A6").Value = "E.Jubilado"
$ws3.Range("B6").Value = 120
$ws3.Range("C6").Value = 4.525
$ws3.Range("D6").Value = 4.53

$ws3.Range("C7").Value = 4.517241379310345
$ws3.Range("D7").Value = 4.52
